$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 247.5
$ws.Range("I2").Value = 112.53846
$ws.Range("K2").Value = 112.53846
$ws.Range("M2").Value = 0.4615399999999994

$ws.Range("H38").Value = 869
$ws.Range("I38").Value = 183.75
$ws.Range("J38").Value = 3610
$ws.Range("K38").Value = 551.25
$ws.Range("L38").Value = 10830
$ws.Range("M38").Value = -179.25
$ws.Range("N38").Value = -11574

$ws.Range("H41").Value = 114.041664
$ws.Range("I41").Value = 228.77777
$ws.Range("J41").Value = 45.2
$ws.Range("K41").Value = 228.77777
$ws.Range("L41").Value = 45.2
$ws.Range("M41").Value = 211.22223
$ws.Range("N41").Value = -925.2

$ws.Range("H116").Value = 8298.963
$ws.Range("I116").Value = 7843.0625
$ws.Range("J116").Value = 8962.091
$ws.Range("K116").Value = 7843.0625
$ws.Range("L116").Value = 8962.091
$ws.Range("M116").Value = -4401.0625
$ws.Range("N116").Value = -15846.091

$ws.Range("H132").Value = 7061.4053
$ws.Range("I132").Value = 4663.48
$ws.Range("K132").Value = 13990.44
$ws.Range("M132").Value = -11460.44

$ws.Range("H137").Value = 3751.7234
$ws.Range("I137").Value = 2669.7144
$ws.Range("J137").Value = 5346.263
$ws.Range("K137").Value = 8009.1432
$ws.Range("L137").Value = 16038.789
$ws.Range("M137").Value = -5459.1432
$ws.Range("N137").Value = -21138.789

$ws.Range("H138").Value = 9196.25
$ws.Range("J138").Value = 9217.956
$ws.Range("L138").Value = 27653.868
$ws.Range("N138").Value = -37933.868

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1922.8077
$ws.Range("I32").Value = 576.4545000000001
$ws.Range("K32").Value = 576.4545000000001
$ws.Range("M32").Value = -289.4545000000001

$ws.Range("H61").Value = 5616.815
$ws.Range("I61").Value = 5716.227
$ws.Range("K61").Value = 5716.227
$ws.Range("M61").Value = -5504.227

$ws.Range("H74").Value = 4619.625
$ws.Range("I74").Value = 4486.75
$ws.Range("J74").Value = 4752.5
$ws.Range("K74").Value = 4486.75
$ws.Range("L74").Value = 4752.5
$ws.Range("M74").Value = -3612.75
$ws.Range("N74").Value = -6500.5

$ws.Range("H77").Value = 4619.625
$ws.Range("I77").Value = 4486.75
$ws.Range("J77").Value = 4752.5
$ws.Range("K77").Value = 22433.75
$ws.Range("L77").Value = 23762.5
$ws.Range("M77").Value = -18065.75
$ws.Range("N77").Value = -32498.5

$ws.Range("H102").Value = 54086.05
$ws.Range("I102").Value = 62960.35
$ws.Range("J102").Value = 3798.3333
$ws.Range("K102").Value = 62960.35
$ws.Range("L102").Value = 3798.3333
$ws.Range("M102").Value = -61338.35
$ws.Range("N102").Value = -7042.3333

$ws.Range("H122").Value = 4318.2905
$ws.Range("I122").Value = 4133.2593
$ws.Range("J122").Value = 5567.25
$ws.Range("K122").Value = 12399.7779
$ws.Range("L122").Value = 16701.75
$ws.Range("M122").Value = -9949.777899999999
$ws.Range("N122").Value = -21601.75

$ws.Range("H132").Value = 121823.25
$ws.Range("I132").Value = 8962.666999999999
$ws.Range("J132").Value = 189539.6
$ws.Range("K132").Value = 26888.001
$ws.Range("L132").Value = 568618.8
$ws.Range("M132").Value = -24358.001
$ws.Range("N132").Value = -573678.8

$ws.Range("H136").Value = 5616.815
$ws.Range("I136").Value = 5716.227
$ws.Range("K136").Value = 17148.681
$ws.Range("M136").Value = -14598.681

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 22749.5
$ws.Range("I97").Value = 15799.4
$ws.Range("K97").Value = 15799.4
$ws.Range("M97").Value = -14808.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4374.0684
$ws.Range("I31").Value = 959.75
$ws.Range("J31").Value = 5045.738
$ws.Range("K31").Value = 959.75
$ws.Range("L31").Value = 5045.738
$ws.Range("M31").Value = -664.75
$ws.Range("N31").Value = -5635.738

$ws.Range("H34").Value = 4374.0684
$ws.Range("I34").Value = 959.75
$ws.Range("J34").Value = 5045.738
$ws.Range("K34").Value = 959.75
$ws.Range("L34").Value = 5045.738
$ws.Range("M34").Value = -757.75
$ws.Range("N34").Value = -5449.738

$ws.Range("H58").Value = 3013.4849
$ws.Range("I58").Value = 2244.72
$ws.Range("J58").Value = 5415.875
$ws.Range("K58").Value = 2244.72
$ws.Range("L58").Value = 5415.875
$ws.Range("M58").Value = -2041.72
$ws.Range("N58").Value = -5821.875

$ws.Range("H86").Value = 7989.18
$ws.Range("I86").Value = 7033.7095
$ws.Range("J86").Value = 9548.105
$ws.Range("K86").Value = 7033.7095
$ws.Range("L86").Value = 9548.105
$ws.Range("M86").Value = -5910.7095
$ws.Range("N86").Value = -11794.105

$ws.Range("H89").Value = 7989.18
$ws.Range("I89").Value = 7033.7095
$ws.Range("J89").Value = 9548.105
$ws.Range("K89").Value = 35168.5475
$ws.Range("L89").Value = 47740.52499999999
$ws.Range("M89").Value = -29552.5475
$ws.Range("N89").Value = -58972.52499999999

$ws.Range("H122").Value = 2402.9333
$ws.Range("I122").Value = 1844.25
$ws.Range("J122").Value = 3041.4285
$ws.Range("K122").Value = 5532.75
$ws.Range("L122").Value = 9124.2855
$ws.Range("M122").Value = -3082.75
$ws.Range("N122").Value = -14024.2855

$ws.Range("H132").Value = 2272.5642
$ws.Range("I132").Value = 1809.3214
$ws.Range("J132").Value = 3451.7273
$ws.Range("K132").Value = 5427.9642
$ws.Range("L132").Value = 10355.1819
$ws.Range("M132").Value = -2897.9642
$ws.Range("N132").Value = -15415.1819

$ws.Range("H134").Value = 361801.12
$ws.Range("I134").Value = 1273136.1
$ws.Range("J134").Value = 7393.0557
$ws.Range("K134").Value = 3819408.3
$ws.Range("L134").Value = 22179.1671
$ws.Range("M134").Value = -3816873.3
$ws.Range("N134").Value = -27249.1671

$ws.Range("H136").Value = 3013.4849
$ws.Range("I136").Value = 2244.72
$ws.Range("J136").Value = 5415.875
$ws.Range("K136").Value = 6734.16
$ws.Range("L136").Value = 16247.625
$ws.Range("M136").Value = -4184.16
$ws.Range("N136").Value = -21347.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 231349.77
$ws.Range("I68").Value = 1499.125
$ws.Range("J68").Value = 362693
$ws.Range("K68").Value = 4497.375
$ws.Range("L68").Value = 1088079
$ws.Range("M68").Value = -3686.375
$ws.Range("N68").Value = -1089701

$ws.Range("H71").Value = 231349.77
$ws.Range("I71").Value = 1499.125
$ws.Range("J71").Value = 362693
$ws.Range("K71").Value = 13492.125
$ws.Range("L71").Value = 3264237
$ws.Range("M71").Value = -9436.125
$ws.Range("N71").Value = -3272349

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 9289.134
$ws.Range("I113").Value = 5068.8
$ws.Range("J113").Value = 11399.3
$ws.Range("K113").Value = 5068.8
$ws.Range("L113").Value = 11399.3
$ws.Range("M113").Value = -2898.8
$ws.Range("N113").Value = -15739.3

$ws.Range("H132").Value = 5700.4614
$ws.Range("I132").Value = 5670.6
$ws.Range("J132").Value = 5800
$ws.Range("K132").Value = 17011.8
$ws.Range("L132").Value = 17400
$ws.Range("M132").Value = -14481.8
$ws.Range("N132").Value = -22460

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5133.8184
$ws.Range("I7").Value = 4380.385
$ws.Range("J7").Value = 5623.55
$ws.Range("K7").Value = 4380.385
$ws.Range("L7").Value = 5623.55
$ws.Range("M7").Value = -4268.385
$ws.Range("N7").Value = -5847.55

$ws.Range("H55").Value = 2676.0667
$ws.Range("I55").Value = 1374.1111
$ws.Range("J55").Value = 4629
$ws.Range("K55").Value = 1374.1111
$ws.Range("L55").Value = 4629
$ws.Range("M55").Value = -1201.1111
$ws.Range("N55").Value = -4975

$ws.Range("H61").Value = 13741.257
$ws.Range("J61").Value = 18879.75
$ws.Range("L61").Value = 18879.75
$ws.Range("N61").Value = -19283.75

$ws.Range("H100").Value = 2733.9583
$ws.Range("I100").Value = 3971.7
$ws.Range("K100").Value = 3971.7
$ws.Range("M100").Value = -3430.7

$ws.Range("H113").Value = 13741.257
$ws.Range("J113").Value = 18879.75
$ws.Range("L113").Value = 18879.75
$ws.Range("N113").Value = -23219.75

$ws.Range("H122").Value = 4144.5264
$ws.Range("I122").Value = 3156.5334
$ws.Range("K122").Value = 9469.600199999999
$ws.Range("M122").Value = -7019.600199999999

$ws.Range("H126").Value = 5133.8184
$ws.Range("I126").Value = 4380.385
$ws.Range("J126").Value = 5623.55
$ws.Range("K126").Value = 13141.155
$ws.Range("L126").Value = 16870.65
$ws.Range("M126").Value = -10671.155
$ws.Range("N126").Value = -21810.65

$ws.Range("H132").Value = 2625.2058
$ws.Range("I132").Value = 2605.6072
$ws.Range("J132").Value = 2716.6667
$ws.Range("K132").Value = 7816.821599999999
$ws.Range("L132").Value = 8150.000100000001
$ws.Range("M132").Value = -5286.821599999999
$ws.Range("N132").Value = -13210.0001

$ws.Range("H136").Value = 25643264
$ws.Range("I136").Value = 1531.1482
$ws.Range("K136").Value = 4593.444600000001
$ws.Range("M136").Value = -2043.444600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 22220.4
$ws.Range("J54").Value = 22220.4
$ws.Range("L54").Value = 22220.4
$ws.Range("N54").Value = -23260.4

$ws.Range("H81").Value = 4766.6665
$ws.Range("J81").Value = 4766.6665
$ws.Range("L81").Value = 9533.333000000001
$ws.Range("N81").Value = -11655.333

$ws.Range("H84").Value = 4766.6665
$ws.Range("J84").Value = 4766.6665
$ws.Range("L84").Value = 47666.665
$ws.Range("N84").Value = -58274.665

$ws.Range("H100").Value = 1439.8889
$ws.Range("I100").Value = 1432.375
$ws.Range("K100").Value = 2864.75
$ws.Range("M100").Value = -2323.75

$ws.Range("H136").Value = 1974.6792
$ws.Range("I136").Value = 1374.2307
$ws.Range("K136").Value = 4122.2307
$ws.Range("M136").Value = -1572.6921
